$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the chi-square calculation formulas for rows 3-5 (Blue, Gold, Red)
foreach ($r in 3..5) {
    $ws.Range("D$r").Formula = "=SUM(B$r,C$r)"
    $ws.Range("E$r").Formula = "=D$r/B7"
    $ws.Range("F$r").Formula = "=D$r/C7"
    $ws.Range("G$r").Formula = "=B$r-E$r"
    $ws.Range("H$r").Formula = "=C$r-F$r"
    $ws.Range("I$r").Formula = "=G$r^2"
    $ws.Range("J$r").Formula = "=H$r^2"
    $ws.Range("K$r").Formula = "=I$r/E$r"
    $ws.Range("L$r").Formula = "=J$r/F$r"
}

# Update the totals row to cover the actual data range (rows 3:5) instead of the empty row 6
$ws.Range("B7").Formula = "=SUM(B3:B5)"
$ws.Range("C7").Formula = "=SUM(C3:C5)"
$ws.Range("K7").Formula = "=SUM(K3:K5)"
$ws.Range("L7").Formula = "=SUM(L3:L5)"
$ws.Range("M7").Formula = "=SUM(K7:L7)"

# Add the chi-square result and degrees of freedom values
$ws.Range("D11").Value = 137476.37
$ws.Range("D12").Value = 2

# Update the current selection
$ws.Range("F12").Select()
